$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row cells:
#   "<Name>_old" -> "<Name>_FV2404"
#   "<Name>_new" -> "<Name>_FV2410"
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2404")
    } elseif ($val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2410")
    }
}

# Turn the header + data range into a native Excel table ("Table1") so the
# header row exposes filter/sort UI, matching the columns we just renamed.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), 0, 1)
$tbl.Name = "Table1"

# Freeze the header row (split/freeze below row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
